$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision process concluded without a selection for a movie on Friday, and therefore, no movie will be acquired.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded successfully.`n"
$ws.Range("D3").Value = "both_movies, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D4").Value = "both_movies, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D5").Value = "no_decision, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie for Friday, resulting in no decision being made.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie.`"`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision process did not result in a selection for Friday's movie, as the committee reached no agreement.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been successfully recorded, and `"Barbie`" will be acquired for Friday’s screening.`n"
$ws.Range("D10").Value = "Barbie_was_selected, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision was recorded as `"no decision.`"`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: I have successfully recorded the decision to acquire the rights for both movies.`n"
$ws.Range("D13").Value = "both_movies, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday has not been made, as the committee did not reach an agreement.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was left unresolved, so no selection has been made.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for Friday's showing.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be resolved, resulting in no selection being made.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: I've successfully recorded the decision to acquire rights for both movies.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: There was no decision made regarding which movie to show on Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("D24").Value = "both_movies, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday resulted in no agreement.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be selected for Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been successfully made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie remains unresolved.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: I have successfully acquired the rights to both movies for the Friday showing.`n"
$ws.Range("D31").Value = "both_movies, "
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: It seems there was no decision reached regarding which movie to show on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: I have recorded the decision about the movie acquisition as `"no decision`".`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The committee did not reach a decision about what movie to show on Friday, so no action is needed regarding movie acquisition.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("D37").Value = "both_movies, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has been classified as a `"no decision.`" No further action regarding the acquisition of movie rights is needed at this time.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: No decision about Friday’s movie was made.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been successfully recorded as `"no decision,`" indicating that the committee could not reach an agreement on which movie to show on Friday.`n"
